# Updates cryptos list data (Price and Volume(1h) columns, plus a few
# Coin/Link cells for rows that were re-ranked) to match the refreshed
# scrape results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.163.09'
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").Value = '''1.568.84'
$ws.Range("E3").Value = '  +1.08%  '

$ws.Range("D4").Value = '''1.01'
$ws.Range("E4").Value = '  +0.89%  '

$ws.Range("D5").Value = '''211.62'
$ws.Range("E5").Value = '  +2.58%  '

$ws.Range("E6").Value = '  +0.81%  '

$ws.Range("E7").Value = '  +1.13%  '

$ws.Range("D8").Value = '''21.96'
$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = '''0.248'
$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("D10").Value = '''0.0597'
$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("D11").Value = '''0.0864'
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").Value = '''1.794.37'
$ws.Range("E12").Value = '  +1.23%  '

$ws.Range("D13").Value = '''1.570.47'
$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("D14").Value = '''3.77'
$ws.Range("E14").Value = '  +0.66%  '

$ws.Range("D15").Value = '''0.519'
$ws.Range("E15").Value = '  +0.03%  '

$ws.Range("D16").Value = '''27.176.05'
$ws.Range("E16").Value = '  +0.99%  '

$ws.Range("D17").Value = '''62.19'
$ws.Range("E17").Value = '  +0.92%  '

$ws.Range("E18").Value = '  -1.18%  '

$ws.Range("D19").Value = '''215.59'
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("E21").Value = '  +1.08%  '

$ws.Range("D22").Value = '''4.13'
$ws.Range("E22").Value = '  +1.33%  '

$ws.Range("D23").Value = '''9.18'
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").Value = '''154.53'
$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("D26").Value = '''6.60'
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("D27").Value = '''15.09'
$ws.Range("E27").Value = '  +0.73%  '

$ws.Range("E28").Value = '  +1.40%  '

$ws.Range("E29").Value = '  +0.96%  '

$ws.Range("E30").Value = '  +5.80%  '

$ws.Range("D31").Value = '''0.0472'
$ws.Range("E31").Value = '  +0.73%  '

$ws.Range("D32").Value = '''3.24'
$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("E33").Value = '  +2.84%  '

$ws.Range("D34").Value = '''1.431.80'
$ws.Range("E34").Value = '  +1.45%  '

$ws.Range("D35").Value = '''1.09'
$ws.Range("E35").Value = '  +12.29%  '

$ws.Range("D37").Value = '''2.36'
$ws.Range("E37").Value = '  +2.69%  '

$ws.Range("E38").Value = '  +1.16%  '

$ws.Range("D39").Value = '''0.530'
$ws.Range("E39").Value = '  +0.66%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''5.85'
$ws.Range("E40").Value = '  +3.49%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''0.808'
$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.01'
$ws.Range("E42").Value = '  +1.17%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.36'
$ws.Range("E43").Value = '  +2.21%  '

$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  +0.87%  '

$ws.Range("D45").Value = '''64.45'
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("D47").Value = '''1.710.01'
$ws.Range("E47").Value = '  +1.37%  '

$ws.Range("D48").Value = '''85.84'
$ws.Range("E48").Value = '  -1.52%  '

$ws.Range("D49").Value = '''0.0₆0101'
$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("D50").Value = '''0.0516'
$ws.Range("E50").Value = '  -1.18%  '

$ws.Range("D51").Value = '''0.0961'
$ws.Range("E51").Value = '  +0.27%  '
